$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-18 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-19 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("236÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "856÷8=", 2) | Out-Null
$d.Content.Find.Execute("537÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "673÷8=", 2) | Out-Null
$d.Content.Find.Execute("413÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "142÷4=", 2) | Out-Null
$d.Content.Find.Execute("238÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "821÷5=", 2) | Out-Null
$d.Content.Find.Execute("968÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "719÷3=", 2) | Out-Null
$d.Content.Find.Execute("621÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "655÷9=", 2) | Out-Null
$d.Content.Find.Execute("154÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "947÷3=", 2) | Out-Null
$d.Content.Find.Execute("376÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "680÷5=", 2) | Out-Null
$d.Content.Find.Execute("401÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "603÷6=", 2) | Out-Null
$d.Content.Find.Execute("615÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "628÷3=", 2) | Out-Null
$d.Content.Find.Execute("347÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "935÷6=", 2) | Out-Null
$d.Content.Find.Execute("172÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "161÷9=", 2) | Out-Null
$d.Content.Find.Execute("343÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "510÷9=", 2) | Out-Null
$d.Content.Find.Execute("695÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "765÷5=", 2) | Out-Null
$d.Content.Find.Execute("156÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "318÷6=", 2) | Out-Null
$d.Content.Find.Execute("888÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "948÷2=", 2) | Out-Null
$d.Content.Find.Execute("182÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "270÷6=", 2) | Out-Null
$d.Content.Find.Execute("673÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "194÷3=", 2) | Out-Null
$d.Content.Find.Execute("106÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "731÷5=", 2) | Out-Null
$d.Content.Find.Execute("975÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "625÷8=", 2) | Out-Null
$d.Content.Find.Execute("406÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "812÷2=", 2) | Out-Null
$d.Content.Find.Execute("420÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "389÷6=", 2) | Out-Null
$d.Content.Find.Execute("408÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "394÷5=", 2) | Out-Null
$d.Content.Find.Execute("114÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "729÷9=", 2) | Out-Null
$d.Content.Find.Execute("193÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "899÷4=", 2) | Out-Null
